# Chapter 3: Comparing Decomposition Implementations - Supervisor suggestions
# Remove unnecessarily detailed text. Reuse important information from the
# Bandwidth section in the Speedup section and remove the rest.
#
# Concretely: the worksheet had four "*-bandwidth" columns (C:F) followed by
# four "*-speedup" columns (G:J). The bandwidth columns are removed entirely
# (their data is no longer needed) and the speedup columns shift left to take
# their place (C:F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the bandwidth columns (cm-bandwidth, icm8-bandwidth, icm16-bandwidth,
# icm32-bandwidth) and delete them, shifting the speedup columns left.
$deleteRange = $ws.Range("C1:F1")
$deleteRange.EntireColumn.Select() | Out-Null
$deleteRange.EntireColumn.Delete() | Out-Null

# Leave the resulting columns (now C:F, the speedup columns) selected, as
# would naturally remain after deleting the preceding columns in Excel.
$ws.Range("C1:F1048576").Select() | Out-Null
